$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Table layout (1-indexed rows/cells):
#  Row 1:  header (Model, K, AIC, dAIC, R2)
#  Row 2:  "Group vs. solo foraging" section label
#  Row 3:  Abiotic + Biotic | 10.00 | 3,913.36 | 0.00 | 0.00  -> R2 0.00 -> 0.11
#  Row 4:  Abiotic           |  7.00 | 3,921.55 | 8.19 | 0.00  -> R2 0.00 -> 0.10
#  Row 5:  Biotic             |  5.00 | 4,094.04 | 180.68 | 0.00 -> R2 0.00 -> 0.11
#  Row 6:  Null
#  Row 7:  "Group size" section label
#  Row 8-11: Group size models (unchanged)
#  Row 12: "Group leadership" section label
#  Row 13: Status            | 4.00 | 2,150.11 | 0.00 | 0.01   -> R2 0.01 -> 0.05
#  Row 14: Null (unchanged)
#  Row 15: Status + Injury   | 6.00 | 2,152.89 | 2.78 | 0.01   -> R2 0.01 -> 0.05
#  Row 16: Injury            | 4.00 | 2,154.76 | 4.65 | 0.01   -> R2 0.01 -> 0.05

# Each data row's R2 value is the 5th (last) cell.
$updates = @(
    @{ Row = 3;  New = "0.11" },
    @{ Row = 4;  New = "0.10" },
    @{ Row = 5;  New = "0.11" },
    @{ Row = 13; New = "0.05" },
    @{ Row = 15; New = "0.05" },
    @{ Row = 16; New = "0.05" }
)

foreach ($u in $updates) {
    $cell = $t.Rows.Item($u.Row).Cells.Item(5)
    $cell.Range.Text = $u.New
}
